$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("A2").Value = 4
$ws.Range("A5").Value = 2

# Add new row of data
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 6

# Update selection to A3
$ws.Range("A3").Select()
